$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 and A3/B3 per diff
$ws.Range("B2").Value = 40
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 17

# Delete rows 4 and 5 (shift cells up) since they're removed entirely
$ws.Range("A4:B5").Delete()
